$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.096.96'
$ws.Cells.Item(2, 5).Value = '  +0.61%  '

$ws.Cells.Item(3, 4).Value = '3.147.48'
$ws.Cells.Item(3, 5).Value = '  +0.62%  '

$ws.Cells.Item(4, 5).Value = '  +0.06%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '592.30'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.58%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '146.27'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.45%  '

$ws.Cells.Item(7, 5).Value = '  +0.03%  '

$ws.Cells.Item(8, 4).Value = '3.141.08'
$ws.Cells.Item(8, 5).Value = '  +0.62%  '

$ws.Cells.Item(9, 5).Value = '  -0.44%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.163'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +1.44%  '

$ws.Cells.Item(11, 5).Value = '  +2.78%  '

$ws.Cells.Item(12, 5).Value = '  -2.15%  '

$ws.Cells.Item(13, 5).Value = '  -2.14%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '37.22'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.20%  '

$ws.Cells.Item(15, 4).Value = '3.668.42'
$ws.Cells.Item(15, 5).Value = '  +0.44%  '

$ws.Cells.Item(16, 5).Value = '  -1.26%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '7.30'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +1.86%  '

$ws.Cells.Item(18, 4).Value = '63.946.11'
$ws.Cells.Item(18, 5).Value = '  +0.51%  '

$ws.Cells.Item(19, 4).Value = '3.146.31'
$ws.Cells.Item(19, 5).Value = '  +0.53%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '467.62'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.59%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '14.34'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.29%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.733'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.10%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '7.51'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.49%  '

$ws.Cells.Item(24, 5).Value = '  -1.95%  '

$ws.Cells.Item(25, 2).Value = 'Litecoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '81.37'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -1.08%  '

$ws.Cells.Item(26, 2).Value = 'Fetch.AI'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.31'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +5.58%  '

$ws.Cells.Item(27, 5).Value = '  +0.11%  '

$ws.Cells.Item(28, 5).Value = '  +7.95%  '

$ws.Cells.Item(29, 5).Value = '  +0.77%  '

$ws.Cells.Item(30, 5).Value = '  +0.08%  '

$ws.Cells.Item(31, 5).Value = '  +6.91%  '

$ws.Cells.Item(32, 5).Value = '  +0.13%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '27.63'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +1.84%  '

$ws.Cells.Item(34, 5).Value = '  +1.50%  '

$ws.Cells.Item(35, 4).Value = '0.0₃0834'
$ws.Cells.Item(35, 5).Value = '  -4.90%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.07'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +1.72%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '6.17'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +0.79%  '

$ws.Cells.Item(38, 5).Value = '  -2.68%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '3.24'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -5.56%  '

$ws.Cells.Item(40, 2).Value = 'Bittensor'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '459.82'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +1.10%  '

$ws.Cells.Item(41, 2).Value = 'OKB'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '51.29'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +0.47%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '9.19'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +5.19%  '

$ws.Cells.Item(43, 5).Value = '  +4.71%  '

$ws.Cells.Item(44, 5).Value = '  -0.41%  '

$ws.Cells.Item(45, 4).Value = '2.925.36'
$ws.Cells.Item(45, 5).Value = '  +0.74%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '39.49'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +10.45%  '

$ws.Cells.Item(47, 5).Value = '  -2.33%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '131.98'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +3.01%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '2.25'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +2.88%  '
